# Regenerate the "K" column (col G) of the save_data sheet.
# Source data was regenerated upstream (std/mean recalculated and s_vals
# recalculated/rewritten), which changed the strike-count ("K") values
# recorded per row. Apply the refreshed K values to the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 2
    6  = 1
    7  = 2
    8  = 1
    9  = 3
    10 = 1
    11 = 2
    12 = 0
    13 = 1
    14 = 3
    15 = 1
    16 = 1
    17 = 0
    18 = 2
    19 = 0
    20 = 2
    21 = 2
    22 = 2
    23 = 1
    24 = 0
    25 = 1
    26 = 0
    27 = 1
    28 = 1
    29 = 1
    31 = 0
    32 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}

Write-Host "Updated K column (G2:G32) for $($newK.Count) rows"
